# Apply the "Add files via upload" edit to the Entradas / Saídas sheets:
#  - Entradas: append 8 new sponsorship rows (rows 4-10)
#  - Saídas: fix the mis-typed date in B3 and append 4 new supplier rows (rows 4-7)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Entradas " - add rows 4..10
# ---------------------------------------------------------------------------
$wsEntradas = $wb.Worksheets.Item("Entradas ")

$entradasRows = @(
    @(44565, 5501, 2),
    @(44778, 5502, 3),
    @(44779, 5503, 4),
    @(44780, 5504, 5),
    @(44781, 5505, 6),
    @(44570, 5506, 7),
    @(44571, 5507, 8)
)

# Row 3 already carries the exact formatting (date style + centered number
# style) that the new rows must reuse, so copy it down instead of rebuilding
# number formats from scratch (that would register new custom numFmts).
$wsEntradas.Range("A3:H3").Copy() | Out-Null
$wsEntradas.Range("A4:H10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$r = 4
foreach ($row in $entradasRows) {
    $wsEntradas.Cells.Item($r, 1).Value = "Patrocínio"
    $wsEntradas.Cells.Item($r, 2).Value = $row[0]
    $wsEntradas.Cells.Item($r, 3).Value = $row[1]
    $wsEntradas.Cells.Item($r, 4).Value = "Ambev"
    $wsEntradas.Cells.Item($r, 5).Value = $row[2]
    $wsEntradas.Cells.Item($r, 6).Value = "-"
    $wsEntradas.Cells.Item($r, 7).Value = "-"
    $wsEntradas.Cells.Item($r, 8).Value = "-"
    $r = $r + 1
}

$wsEntradas.Range("D13").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Saídas" - fix B3 date, add rows 4..7
# ---------------------------------------------------------------------------
$wsSaidas = $wb.Worksheets.Item("Saídas")

# B3 held the literal text "02/012020" - replace with the real date (2022-01-02).
# The cell already carries the right numFmt (d-mmm, style 4), so just overwrite
# the value - no format change needed.
$wsSaidas.Cells.Item(3, 2).Value = 44563

# Row 3 has the number formats the new rows need (date style on B, plain on
# the rest), so copy it down before filling in the actual values.
$wsSaidas.Range("A3:E3").Copy() | Out-Null
$wsSaidas.Range("A4:E7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$saidasRows = @(
    @(44564, 1564.52),
    @(44565, 1565.52),
    @(44778, 1566.52),
    @(44779, 1567.52)
)

$r = 4
foreach ($row in $saidasRows) {
    $wsSaidas.Cells.Item($r, 1).Value = "Pagamento"
    $wsSaidas.Cells.Item($r, 2).Value = $row[0]
    $wsSaidas.Cells.Item($r, 3).Value = $row[1]
    $wsSaidas.Cells.Item($r, 4).Value = "João Flávio"
    $wsSaidas.Cells.Item($r, 5).Value = "Fornecedor de destilados"
    $r = $r + 1
}

$wsSaidas.Range("B3").Select() | Out-Null
